$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the missing hours value for "Cody Note" (row 6) in column B (Week 1)
$ws.Range("B6").Value = 15

# Update the active selection to B7 (also clears any custom topLeftCell scroll position)
$ws.Range("B7").Select()
